$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 2.590158771791494
$ws.Range("C2").Value = 0.6531155350381255
$ws.Range("D2").Value = 0.04807882700767863
$ws.Range("E2").Value = 0.1195447155408011
$ws.Range("F2").Value = 3.384970243692806
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.2274430023337857
$ws.Range("N2").Value = 1.807420711823944

$ws.Range("B3").Value = 2.437359848274127
$ws.Range("C3").Value = 0.6069069128337787
$ws.Range("D3").Value = 0.04793712971785169
$ws.Range("E3").Value = 0.1177158772568525
$ws.Range("F3").Value = 3.328819371221186
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.2223673216307844
$ws.Range("N3").Value = 1.822854897791373

$ws.Range("B4").Value = 2.345226233922517
$ws.Range("C4").Value = 0.5789541211845517
$ws.Range("D4").Value = 0.04786182353971924
$ws.Range("E4").Value = 0.116659652562582
$ws.Range("F4").Value = 3.296509297508322
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.219397220266373
$ws.Range("N4").Value = 1.832999409062673

$ws.Range("B5").Value = 2.308100492833603
$ws.Range("C5").Value = 0.5676671481167546
$ws.Range("D5").Value = 0.04783410073699201
$ws.Range("E5").Value = 0.1162459481881797
$ws.Range("F5").Value = 3.283883502477281
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.218223428447395
$ws.Range("N5").Value = 1.837300489266624

$ws.Range("B6").Value = 2.301960999330333
$ws.Range("C6").Value = 0.5657991948310155
$ws.Range("D6").Value = 0.04782967717854092
$ws.Range("E6").Value = 0.116178260947482
$ws.Range("F6").Value = 3.281819546344792
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.2180307210998933
$ws.Range("N6").Value = 1.83802474926653

$ws.Range("B7").Value = 2.344723850566822
$ws.Range("C7").Value = 0.5788014818782585
$ws.Range("D7").Value = 0.0478614376230091
$ws.Range("E7").Value = 0.1166540055885754
$ws.Range("F7").Value = 3.296336837125949
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.2193812424567412
$ws.Range("N7").Value = 1.83305673939703

$ws.Range("B8").Value = 2.537121150666167
$ws.Range("C8").Value = 0.6370948311719076
$ws.Range("D8").Value = 0.04802755628642075
$ws.Range("E8").Value = 0.1189002569276987
$ws.Range("F8").Value = 3.365157235577044
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.2256623641175253
$ws.Range("N8").Value = 1.812603304945881

$ws.Range("B9").Value = 2.928008787594251
$ws.Range("C9").Value = 0.754809882341533
$ws.Range("D9").Value = 0.04844512876813667
$ws.Range("E9").Value = 0.1238371838019923
$ws.Range("F9").Value = 3.517501306237563
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.2391540164819617
$ws.Range("N9").Value = 1.77783049327212

$ws.Range("B10").Value = 3.223828490996254
$ws.Range("C10").Value = 0.8434814779427597
$ws.Range("D10").Value = 0.04880664150051572
$ws.Range("E10").Value = 0.1277933445843864
$ws.Range("F10").Value = 3.64031728116305
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.249801667486409
$ws.Range("N10").Value = 1.755587190251418

$ws.Range("B11").Value = 3.360357790341993
$ws.Range("C11").Value = 0.8843206251112292
$ws.Range("D11").Value = 0.04898273104867812
$ws.Range("E11").Value = 0.1296656243145478
$ws.Range("F11").Value = 3.698618511605048
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.2548095989843659
$ws.Range("N11").Value = 1.746196872268186

$ws.Range("B12").Value = 3.41234550929795
$ws.Range("C12").Value = 0.8998594549816517
$ws.Range("D12").Value = 0.04905106291874972
$ws.Range("E12").Value = 0.1303851281592259
$ws.Range("F12").Value = 3.721050406410228
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.2567299308575599
$ws.Range("N12").Value = 1.742746698872224

$ws.Range("B13").Value = 3.401136172465044
$ws.Range("C13").Value = 0.8965095736126614
$ws.Range("D13").Value = 0.04903627337458261
$ws.Range("E13").Value = 0.1302297014287035
$ws.Range("F13").Value = 3.716203452947582
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.2563152833179885
$ws.Range("N13").Value = 1.743485036492572

$ws.Range("B14").Value = 3.36462906796686
$ws.Range("C14").Value = 0.885597522161504
$ws.Range("D14").Value = 0.04898831980093732
$ws.Range("E14").Value = 0.1297246072281766
$ws.Range("F14").Value = 3.700456864757314
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.25496710427457
$ws.Range("N14").Value = 1.745910899461364

$ws.Range("B15").Value = 3.342304951657809
$ws.Range("C15").Value = 0.8789232560894789
$ws.Range("D15").Value = 0.04895916114290699
$ws.Range("E15").Value = 0.1294165936188278
$ws.Range("F15").Value = 3.690857933776982
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.2541444324836561
$ws.Range("N15").Value = 1.747410610864577

$ws.Range("B16").Value = 3.21494585672059
$ws.Range("C16").Value = 0.8408227932669092
$ws.Range("D16").Value = 0.04879536583182897
$ws.Range("E16").Value = 0.1276724537439691
$ws.Range("F16").Value = 3.636556498965888
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.2494777196239681
$ws.Range("N16").Value = 1.756215610390342

$ws.Range("B17").Value = 3.137320745216982
$ws.Range("C17").Value = 0.8175791947599009
$ws.Range("D17").Value = 0.04869784666452404
$ws.Range("E17").Value = 0.1266211260672385
$ws.Range("F17").Value = 3.603870524313919
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.2466571365859664
$ws.Range("N17").Value = 1.761804434660206

$ws.Range("B18").Value = 3.09285694210115
$ws.Range("C18").Value = 0.8042572348392127
$ws.Range("D18").Value = 0.0486428527444609
$ws.Range("E18").Value = 0.1260232594105766
$ws.Range("F18").Value = 3.585298998594396
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.245050253566319
$ws.Range("N18").Value = 1.765087465476526

$ws.Range("B19").Value = 3.077833724528375
$ws.Range("C19").Value = 0.7997546924987091
$ws.Range("D19").Value = 0.04862442177084603
$ws.Range("E19").Value = 0.1258220028051191
$ws.Range("F19").Value = 3.579050113020145
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.2445088331995464
$ws.Range("N19").Value = 1.766210775429542

$ws.Range("B20").Value = 3.145564978607865
$ws.Range("C20").Value = 0.8200486238513349
$ws.Range("D20").Value = 0.04870811442526346
$ws.Range("E20").Value = 0.1267323343815967
$ws.Range("F20").Value = 3.6073263101093
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.246955792035493
$ws.Range("N20").Value = 1.761202399370774

$ws.Range("B21").Value = 3.37534425674977
$ws.Range("C21").Value = 0.8888006356021947
$ws.Range("D21").Value = 0.04900236030996652
$ws.Range("E21").Value = 0.1298726797555076
$ws.Range("F21").Value = 3.705072359926191
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.2553624449436853
$ws.Range("N21").Value = 1.745195486902489

$ws.Range("B22").Value = 3.527193906445632
$ws.Range("C22").Value = 0.9341658166571847
$ws.Range("D22").Value = 0.04920427762529656
$ws.Range("E22").Value = 0.1319863850729206
$ws.Range("F22").Value = 3.771023564472699
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.2609963531017456
$ws.Range("N22").Value = 1.735350837760734

$ws.Range("B23").Value = 3.445993796063476
$ws.Range("C23").Value = 0.9099134716588537
$ws.Range("D23").Value = 0.04909563855170518
$ws.Range("E23").Value = 0.1308526266097587
$ws.Range("F23").Value = 3.735633279034431
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.257976543685345
$ws.Range("N23").Value = 1.740548336851745

$ws.Range("B24").Value = 3.141837250553465
$ws.Range("C24").Value = 0.8189320670904863
$ws.Range("D24").Value = 0.04870346903217282
$ws.Range("E24").Value = 0.126682036688841
$ws.Range("F24").Value = 3.605763264225175
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.2468207241193312
$ws.Range("N24").Value = 1.761474361749904

$ws.Range("B25").Value = 2.820772391448429
$ws.Range("C25").Value = 0.7225891535513256
$ws.Range("D25").Value = 0.04832246564163256
$ws.Range("E25").Value = 0.1224441832956131
$ws.Range("F25").Value = 3.474397201050692
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.2353764591542244
$ws.Range("N25").Value = 1.786660944662032
